$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range('E13').Style = "Percent"
$ws.Range('E13').NumberFormat = "0.0%"
$nf = $ws.Range('E13').NumberFormat
Write-Output "E13 NumberFormat=$nf"
